# PauseTexts.xlsx - "Forgot to commit earlier. Now starting development on
# Input handeling." follow-up edit.
#
# Summary of the change being applied:
#   - A new column is inserted before column E (i.e. a fresh column D),
#     pushing the old D/E/F/G/H/I content one column to the right.
#   - The new D3 cell gets a "followup" header.
#   - D10 (under the new "followup" column, on the row that already carries
#     the highlighted A10/C10 cells) is formatted but left empty.
#   - A brand new row 14 is added with two more text entries.
#   - The active selection is moved to F12 (the relocated "Imperfection"
#     highlighted cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; everything from the old D onward (D,E,F,G,H,I)
# shifts right by one (new E,F,G,H,I,J), carrying its values/styles/widths.
$ws.Columns("D").Insert()

# New header text in the freshly inserted column.
$ws.Range("D3").Value = "followup"

# D10 picked up C10's highlighted style from the column insert; clear its
# fill back off (explicit "no fill") while keeping it as its own distinct
# format, matching the dedicated format added for this cell.
$ws.Range("D10").Interior.ColorIndex = -4142

# Two brand-new cells making up the extra row at the bottom of the table.
$ws.Range("C14").Value = "After pain comes satisfaction"
$ws.Range("D14").Value = "Always."

# Leave the selection on the relocated highlighted cell.
$ws.Range("F12").Select() | Out-Null
